$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("I1").Value = "Electric Choice Id"
$ws.Range("J1").Value = "Rate Code"
$ws.Range("L1").Value = "Usage"

# Remove the Gas columns (M, N, O) entirely - shifts everything left
$ws.Range("M1:O2").EntireColumn.Delete()

# Clear the Electric Supplier value in row 2 (was "N/A")
$ws.Range("K2").ClearContents()
